$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Critical Minutes / Good Roaming Calculation for row 3 and totals row 4
$ws.Range("C3").Value = 1539
$ws.Range("D3").Value = 85.59999999999999
$ws.Range("C4").Value = 1539

# Clear the Driver Vintage date in E12
$ws.Range("E12").ClearContents()

# Update Total Samples for row 14
$ws.Range("B14").Value = 265400
